$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 400  # H12: 450 -> 400
$ws.Cells.Item(12, 9).Value = 0  # I12: 500 -> 0
$ws.Cells.Item(12, 11).Value = 0  # K12: 500 -> 0
$ws.Cells.Item(12, 13).ClearContents()  # M12: -330 -> (removed)
$ws.Cells.Item(18, 8).Value = 2333.3333  # H18: 2000 -> 2333.3333
$ws.Cells.Item(18, 9).Value = 2333.3333  # I18: 2000 -> 2333.3333
$ws.Cells.Item(18, 11).Value = 2333.3333  # K18: 2000 -> 2333.3333
$ws.Cells.Item(18, 13).Value = -2049.3333  # M18: -1716 -> -2049.3333
$ws.Cells.Item(32, 8).Value = 8500  # H32: 8136.364 -> 8500
$ws.Cells.Item(32, 10).Value = 9125  # J32: 8611.111000000001 -> 9125
$ws.Cells.Item(32, 12).Value = 9125  # L32: 8611.111000000001 -> 9125
$ws.Cells.Item(32, 14).Value = -9777  # N32: -9263.111000000001 -> -9777
$ws.Cells.Item(40, 8).Value = 8734.385  # H40: 8089 -> 8734.385
$ws.Cells.Item(40, 9).Value = 3750.5  # I40: 2979.8 -> 3750.5
$ws.Cells.Item(40, 10).Value = 9149.708000000001  # J40: 9199.695 -> 9149.708000000001
$ws.Cells.Item(40, 11).Value = 3750.5  # K40: 2979.8 -> 3750.5
$ws.Cells.Item(40, 12).Value = 9149.708000000001  # L40: 9199.695 -> 9149.708000000001
$ws.Cells.Item(40, 13).Value = -3575.5  # M40: -2804.8 -> -3575.5
$ws.Cells.Item(40, 14).Value = -9499.708000000001  # N40: -9549.695 -> -9499.708000000001
$ws.Cells.Item(64, 8).Value = 3995  # H64: 3999 -> 3995
$ws.Cells.Item(64, 9).Value = 3995  # I64: 3999 -> 3995
$ws.Cells.Item(64, 11).Value = 3995  # K64: 3999 -> 3995
$ws.Cells.Item(64, 13).Value = -3747  # M64: -3751 -> -3747
$ws.Cells.Item(67, 8).Value = 3995  # H67: 3999 -> 3995
$ws.Cells.Item(67, 9).Value = 3995  # I67: 3999 -> 3995
$ws.Cells.Item(67, 11).Value = 3995  # K67: 3999 -> 3995
$ws.Cells.Item(67, 13).Value = -3137  # M67: -3141 -> -3137
$ws.Cells.Item(86, 8).Value = 3000  # H86: 0 -> 3000
$ws.Cells.Item(86, 10).Value = 3000  # J86: 0 -> 3000
$ws.Cells.Item(86, 12).Value = 3000  # L86: 0 -> 3000
$ws.Cells.Item(86, 14).Value = -5246  # N86: None -> -5246
$ws.Cells.Item(89, 8).Value = 3000  # H89: 0 -> 3000
$ws.Cells.Item(89, 10).Value = 3000  # J89: 0 -> 3000
$ws.Cells.Item(89, 12).Value = 15000  # L89: 0 -> 15000
$ws.Cells.Item(89, 14).Value = -26232  # N89: None -> -26232
$ws.Cells.Item(137, 8).Value = 1898.4  # H137: 1792.7368 -> 1898.4
$ws.Cells.Item(137, 9).Value = 1855.9  # I137: 1749.6923 -> 1855.9
$ws.Cells.Item(137, 10).Value = 1983.4  # J137: 1886 -> 1983.4
$ws.Cells.Item(137, 11).Value = 5567.700000000001  # K137: 5249.0769 -> 5567.700000000001
$ws.Cells.Item(137, 12).Value = 5950.200000000001  # L137: 5658 -> 5950.200000000001
$ws.Cells.Item(137, 13).Value = -3017.700000000001  # M137: -2699.0769 -> -3017.700000000001
$ws.Cells.Item(137, 14).Value = -11050.2  # N137: -10758 -> -11050.2
$ws.Cells.Item(138, 8).Value = 3100.3572  # H138: 3031.3462 -> 3100.3572
$ws.Cells.Item(138, 10).Value = 3998.842  # J138: 3999 -> 3998.842
$ws.Cells.Item(138, 12).Value = 11996.526  # L138: 11997 -> 11996.526
$ws.Cells.Item(138, 14).Value = -22276.526  # N138: -22277 -> -22276.526
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 300.33334  # H5: 230.2 -> 300.33334
$ws.Cells.Item(5, 9).Value = 300.33334  # I5: 230.2 -> 300.33334
$ws.Cells.Item(5, 11).Value = 300.33334  # K5: 230.2 -> 300.33334
$ws.Cells.Item(5, 13).Value = -188.33334  # M5: -118.2 -> -188.33334
$ws.Cells.Item(32, 8).Value = 2187  # H32: 2127.0881 -> 2187
$ws.Cells.Item(32, 9).Value = 2187  # I32: 2127.0881 -> 2187
$ws.Cells.Item(32, 11).Value = 2187  # K32: 2127.0881 -> 2187
$ws.Cells.Item(32, 13).Value = -1900  # M32: -1840.0881 -> -1900
$ws.Cells.Item(45, 8).Value = 2130.3333  # H45: 2298.2144 -> 2130.3333
$ws.Cells.Item(45, 9).Value = 1845.3  # I45: 1778.5454 -> 1845.3
$ws.Cells.Item(45, 10).Value = 3555.5  # J45: 4203.6665 -> 3555.5
$ws.Cells.Item(45, 11).Value = 1845.3  # K45: 1778.5454 -> 1845.3
$ws.Cells.Item(45, 12).Value = 3555.5  # L45: 4203.6665 -> 3555.5
$ws.Cells.Item(45, 13).Value = -1468.3  # M45: -1401.5454 -> -1468.3
$ws.Cells.Item(45, 14).Value = -4309.5  # N45: -4957.6665 -> -4309.5
$ws.Cells.Item(46, 8).Value = 12166.333  # H46: 14166.333 -> 12166.333
$ws.Cells.Item(46, 9).Value = 12166.333  # I46: 14166.333 -> 12166.333
$ws.Cells.Item(46, 11).Value = 12166.333  # K46: 14166.333 -> 12166.333
$ws.Cells.Item(46, 13).Value = -11847.333  # M46: -13847.333 -> -11847.333
$ws.Cells.Item(122, 8).Value = 5904.1113  # H122: 6348 -> 5904.1113
$ws.Cells.Item(122, 9).Value = 6654.8423  # I122: 7096.1665 -> 6654.8423
$ws.Cells.Item(122, 10).Value = 4121.125  # J122: 4424.143 -> 4121.125
$ws.Cells.Item(122, 11).Value = 19964.5269  # K122: 21288.4995 -> 19964.5269
$ws.Cells.Item(122, 12).Value = 12363.375  # L122: 13272.429 -> 12363.375
$ws.Cells.Item(122, 13).Value = -17514.5269  # M122: -18838.4995 -> -17514.5269
$ws.Cells.Item(122, 14).Value = -17263.375  # N122: -18172.429 -> -17263.375
$ws.Cells.Item(132, 8).Value = 3944  # H132: 4999.5 -> 3944
$ws.Cells.Item(132, 9).Value = 3944  # I132: 4999.5 -> 3944
$ws.Cells.Item(132, 11).Value = 11832  # K132: 14998.5 -> 11832
$ws.Cells.Item(132, 13).Value = -9302  # M132: -12468.5 -> -9302
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 300.33334  # H4: 230.2 -> 300.33334
$ws.Cells.Item(4, 9).Value = 300.33334  # I4: 230.2 -> 300.33334
$ws.Cells.Item(4, 11).Value = 300.33334  # K4: 230.2 -> 300.33334
$ws.Cells.Item(4, 13).Value = -185.33334  # M4: -115.2 -> -185.33334
$ws.Cells.Item(22, 8).Value = 5556179.5  # H22: 2778664.8 -> 5556179.5
$ws.Cells.Item(22, 9).Value = 5556179.5  # I22: 3704786.2 -> 5556179.5
$ws.Cells.Item(22, 10).Value = 0  # J22: 300 -> 0
$ws.Cells.Item(22, 11).Value = 5556179.5  # K22: 3704786.2 -> 5556179.5
$ws.Cells.Item(22, 12).Value = 0  # L22: 300 -> 0
$ws.Cells.Item(22, 13).Value = -5556006.5  # M22: -3704613.2 -> -5556006.5
$ws.Cells.Item(22, 14).ClearContents()  # N22: -646 -> (removed)
$ws.Cells.Item(82, 8).Value = 24215.2  # H82: 26461.334 -> 24215.2
$ws.Cells.Item(82, 9).Value = 10358.667  # I82: 11630.4 -> 10358.667
$ws.Cells.Item(82, 11).Value = 10358.667  # K82: 11630.4 -> 10358.667
$ws.Cells.Item(82, 13).Value = -9975.666999999999  # M82: -11247.4 -> -9975.666999999999
$ws.Cells.Item(85, 8).Value = 24215.2  # H85: 26461.334 -> 24215.2
$ws.Cells.Item(85, 9).Value = 10358.667  # I85: 11630.4 -> 10358.667
$ws.Cells.Item(85, 11).Value = 10358.667  # K85: 11630.4 -> 10358.667
$ws.Cells.Item(85, 13).Value = -9032.666999999999  # M85: -10304.4 -> -9032.666999999999
$ws.Cells.Item(86, 8).Value = 1166.3334  # H86: 0 -> 1166.3334
$ws.Cells.Item(86, 9).Value = 1250  # I86: 0 -> 1250
$ws.Cells.Item(86, 10).Value = 999  # J86: 0 -> 999
$ws.Cells.Item(86, 11).Value = 1250  # K86: 0 -> 1250
$ws.Cells.Item(86, 12).Value = 999  # L86: 0 -> 999
$ws.Cells.Item(86, 13).Value = -127  # M86: None -> -127
$ws.Cells.Item(86, 14).Value = -3245  # N86: None -> -3245
$ws.Cells.Item(89, 8).Value = 1166.3334  # H89: 0 -> 1166.3334
$ws.Cells.Item(89, 9).Value = 1250  # I89: 0 -> 1250
$ws.Cells.Item(89, 10).Value = 999  # J89: 0 -> 999
$ws.Cells.Item(89, 11).Value = 6250  # K89: 0 -> 6250
$ws.Cells.Item(89, 12).Value = 4995  # L89: 0 -> 4995
$ws.Cells.Item(89, 13).Value = -634  # M89: None -> -634
$ws.Cells.Item(89, 14).Value = -16227  # N89: None -> -16227
$ws.Cells.Item(97, 8).Value = 24000  # H97: 23000 -> 24000
$ws.Cells.Item(97, 9).Value = 28000  # I97: 23000 -> 28000
$ws.Cells.Item(97, 10).Value = 20000  # J97: 0 -> 20000
$ws.Cells.Item(97, 11).Value = 28000  # K97: 23000 -> 28000
$ws.Cells.Item(97, 12).Value = 20000  # L97: 0 -> 20000
$ws.Cells.Item(97, 13).Value = -27009  # M97: -22009 -> -27009
$ws.Cells.Item(97, 14).Value = -21982  # N97: None -> -21982
$ws.Cells.Item(99, 8).Value = 3789.7334  # H99: 3794.9333 -> 3789.7334
$ws.Cells.Item(99, 9).Value = 3989  # I99: 3994.5715 -> 3989
$ws.Cells.Item(99, 11).Value = 3989  # K99: 3994.5715 -> 3989
$ws.Cells.Item(99, 13).Value = -2491  # M99: -2496.5715 -> -2491
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 960  # H22: 1476 -> 960
$ws.Cells.Item(22, 9).Value = 960  # I22: 961.6667 -> 960
$ws.Cells.Item(22, 10).Value = 0  # J22: 2247.5 -> 0
$ws.Cells.Item(22, 11).Value = 960  # K22: 961.6667 -> 960
$ws.Cells.Item(22, 12).Value = 0  # L22: 2247.5 -> 0
$ws.Cells.Item(22, 13).Value = -610  # M22: -611.6667 -> -610
$ws.Cells.Item(22, 14).ClearContents()  # N22: -2947.5 -> (removed)
$ws.Cells.Item(122, 8).Value = 929.1667  # H122: 941.46155 -> 929.1667
$ws.Cells.Item(122, 9).Value = 939  # I122: 950.1111 -> 939
$ws.Cells.Item(122, 10).Value = 899.6667  # J122: 922 -> 899.6667
$ws.Cells.Item(122, 11).Value = 2817  # K122: 2850.3333 -> 2817
$ws.Cells.Item(122, 12).Value = 2699.0001  # L122: 2766 -> 2699.0001
$ws.Cells.Item(122, 13).Value = -367  # M122: -400.3332999999998 -> -367
$ws.Cells.Item(122, 14).Value = -7599.0001  # N122: -7666 -> -7599.0001
$ws.Cells.Item(132, 8).Value = 1385.1666  # H132: 2039 -> 1385.1666
$ws.Cells.Item(132, 9).Value = 937  # I132: 933.5 -> 937
$ws.Cells.Item(132, 10).Value = 1833.3334  # J132: 4250 -> 1833.3334
$ws.Cells.Item(132, 11).Value = 2811  # K132: 2800.5 -> 2811
$ws.Cells.Item(132, 12).Value = 5500.0002  # L132: 12750 -> 5500.0002
$ws.Cells.Item(132, 13).Value = -281  # M132: -270.5 -> -281
$ws.Cells.Item(132, 14).Value = -10560.0002  # N132: -17810 -> -10560.0002
$ws.Cells.Item(134, 8).Value = 2332.0667  # H134: 2353 -> 2332.0667
$ws.Cells.Item(134, 9).Value = 2244.3635  # I134: 2264.9 -> 2244.3635
$ws.Cells.Item(134, 11).Value = 6733.0905  # K134: 6794.700000000001 -> 6733.0905
$ws.Cells.Item(134, 13).Value = -4198.0905  # M134: -4259.700000000001 -> -4198.0905
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 1546.5555  # H18: 2060 -> 1546.5555
$ws.Cells.Item(18, 9).Value = 1546.5555  # I18: 2060 -> 1546.5555
$ws.Cells.Item(18, 11).Value = 4639.666499999999  # K18: 6180 -> 4639.666499999999
$ws.Cells.Item(18, 13).Value = -4470.666499999999  # M18: -6011 -> -4470.666499999999
$ws.Cells.Item(40, 8).Value = 191.6  # H40: 211.22223 -> 191.6
$ws.Cells.Item(40, 9).Value = 76.59999999999999  # I40: 92 -> 76.59999999999999
$ws.Cells.Item(40, 11).Value = 306.4  # K40: 368 -> 306.4
$ws.Cells.Item(40, 13).Value = -237.4  # M40: -299 -> -237.4
$ws.Cells.Item(131, 8).Value = 4664.231  # H131: 0 -> 4664.231
$ws.Cells.Item(131, 9).Value = 635  # I131: 0 -> 635
$ws.Cells.Item(131, 10).Value = 5000  # J131: 0 -> 5000
$ws.Cells.Item(131, 11).Value = 1905  # K131: 0 -> 1905
$ws.Cells.Item(131, 12).Value = 15000  # L131: 0 -> 15000
$ws.Cells.Item(131, 13).Value = 3135  # M131: None -> 3135
$ws.Cells.Item(131, 14).Value = -25080  # N131: None -> -25080
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 3000000  # H52: 2510000 -> 3000000
$ws.Cells.Item(52, 10).Value = 3000000  # J52: 2510000 -> 3000000
$ws.Cells.Item(52, 12).Value = 3000000  # L52: 2510000 -> 3000000
$ws.Cells.Item(52, 14).Value = -3000518  # N52: -2510518 -> -3000518
$ws.Cells.Item(122, 8).Value = 2591.1667  # H122: 2811.5 -> 2591.1667
$ws.Cells.Item(122, 9).Value = 2549.2307  # I122: 2862.0908 -> 2549.2307
$ws.Cells.Item(122, 11).Value = 7647.6921  # K122: 8586.2724 -> 7647.6921
$ws.Cells.Item(122, 13).Value = -5197.6921  # M122: -6136.2724 -> -5197.6921
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7160.125  # H7: 6328 -> 7160.125
$ws.Cells.Item(7, 9).Value = 7046.8335  # I7: 6035 -> 7046.8335
$ws.Cells.Item(7, 11).Value = 7046.8335  # K7: 6035 -> 7046.8335
$ws.Cells.Item(7, 13).Value = -6934.8335  # M7: -5923 -> -6934.8335
$ws.Cells.Item(16, 8).Value = 5683  # H16: 6262.25 -> 5683
$ws.Cells.Item(16, 9).Value = 5683  # I16: 6262.25 -> 5683
$ws.Cells.Item(16, 11).Value = 5683  # K16: 6262.25 -> 5683
$ws.Cells.Item(16, 13).Value = -5513  # M16: -6092.25 -> -5513
$ws.Cells.Item(46, 8).Value = 900  # H46: 850 -> 900
$ws.Cells.Item(46, 9).Value = 800  # I46: 775 -> 800
$ws.Cells.Item(46, 11).Value = 800  # K46: 775 -> 800
$ws.Cells.Item(46, 13).Value = -612  # M46: -587 -> -612
$ws.Cells.Item(81, 8).Value = 47500  # H81: 0 -> 47500
$ws.Cells.Item(81, 10).Value = 47500  # J81: 0 -> 47500
$ws.Cells.Item(81, 12).Value = 47500  # L81: 0 -> 47500
$ws.Cells.Item(81, 14).Value = -49496  # N81: None -> -49496
$ws.Cells.Item(84, 8).Value = 47500  # H84: 0 -> 47500
$ws.Cells.Item(84, 10).Value = 47500  # J84: 0 -> 47500
$ws.Cells.Item(84, 12).Value = 142500  # L84: 0 -> 142500
$ws.Cells.Item(84, 14).Value = -152484  # N84: None -> -152484
$ws.Cells.Item(126, 8).Value = 7160.125  # H126: 6328 -> 7160.125
$ws.Cells.Item(126, 9).Value = 7046.8335  # I126: 6035 -> 7046.8335
$ws.Cells.Item(126, 11).Value = 21140.5005  # K126: 18105 -> 21140.5005
$ws.Cells.Item(126, 13).Value = -18670.5005  # M126: -15635 -> -18670.5005
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 10000  # H14: 9831.666999999999 -> 10000
$ws.Cells.Item(14, 10).Value = 0  # J14: 9747.5 -> 0
$ws.Cells.Item(14, 12).Value = 0  # L14: 9747.5 -> 0
$ws.Cells.Item(14, 14).ClearContents()  # N14: -10083.5 -> (removed)
$ws.Cells.Item(51, 8).Value = 0  # H51: 25035 -> 0
$ws.Cells.Item(51, 9).Value = 0  # I51: 20070 -> 0
$ws.Cells.Item(51, 10).Value = 0  # J51: 30000 -> 0
$ws.Cells.Item(51, 11).Value = 0  # K51: 20070 -> 0
$ws.Cells.Item(51, 12).Value = 0  # L51: 30000 -> 0
$ws.Cells.Item(51, 13).ClearContents()  # M51: -19560 -> (removed)
$ws.Cells.Item(51, 14).ClearContents()  # N51: -31020 -> (removed)
$ws.Cells.Item(52, 8).Value = 2880863.2  # H52: 2523010.5 -> 2880863.2
$ws.Cells.Item(52, 9).Value = 5014510.5  # I52: 4015216.5 -> 5014510.5
$ws.Cells.Item(52, 11).Value = 5014510.5  # K52: 4015216.5 -> 5014510.5
$ws.Cells.Item(52, 13).Value = -5014284.5  # M52: -4014990.5 -> -5014284.5
$ws.Cells.Item(58, 8).Value = 100000000  # H58: 7418.3335 -> 100000000
$ws.Cells.Item(58, 9).Value = 100000000  # I58: 7418.3335 -> 100000000
$ws.Cells.Item(58, 11).Value = 100000000  # K58: 7418.3335 -> 100000000
$ws.Cells.Item(58, 13).Value = -99999692  # M58: -7110.3335 -> -99999692
$ws.Cells.Item(107, 8).Value = 220  # H107: 0 -> 220
$ws.Cells.Item(107, 9).Value = 220  # I107: 0 -> 220
$ws.Cells.Item(107, 11).Value = 660  # K107: 0 -> 660
$ws.Cells.Item(107, 13).Value = 1260  # M107: None -> 1260
$ws.Cells.Item(122, 8).Value = 2561.8  # H122: 1660.04 -> 2561.8
$ws.Cells.Item(122, 9).Value = 1652  # I122: 1413.0869 -> 1652
$ws.Cells.Item(122, 10).Value = 3168.3333  # J122: 4500 -> 3168.3333
$ws.Cells.Item(122, 11).Value = 4956  # K122: 4239.2607 -> 4956
$ws.Cells.Item(122, 12).Value = 9504.999899999999  # L122: 13500 -> 9504.999899999999
$ws.Cells.Item(122, 13).Value = -2506  # M122: -1789.2607 -> -2506
$ws.Cells.Item(122, 14).Value = -14404.9999  # N122: -18400 -> -14404.9999
$ws.Cells.Item(132, 8).Value = 3089.442  # H132: 3337.7026 -> 3089.442
$ws.Cells.Item(132, 9).Value = 2140.7632  # I132: 2249.9375 -> 2140.7632
$ws.Cells.Item(132, 11).Value = 6422.2896  # K132: 6749.8125 -> 6422.2896
$ws.Cells.Item(132, 13).Value = -3892.2896  # M132: -4219.8125 -> -3892.2896
